# Avoid Word's "smart quotes" autocorrect messing with straight quotes we insert.
try { $word.Options.AutoFormatAsYouTypeReplaceQuotes = $false } catch {}
try { $word.Options.AutoFormatReplaceQuotes = $false } catch {}

$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Insert a new "Meta description" paragraph right after the first
#    (Heading1) paragraph: "Play 1 Million Fortunes Megaways free - Review & Guide"
# ---------------------------------------------------------------------------
$metaLabel = "Meta description"
$metaRest  = ": Our review of 1 Million Fortunes Megaways: enjoy up to 1,058,841 ways to win, Asian-themed graphics, and an innovative cascading mechanic. Play for free."

$pHeading = $d.Paragraphs.Item(1)
$pNext = $d.Paragraphs.Item(2)
$pNext.Range.InsertParagraphBefore()

$pMeta = $d.Paragraphs.Item(2)
$pMeta.Style = "Normal"

$rMeta = $pMeta.Range
$rMeta.Collapse(0)
$rMeta.MoveEnd(1, -1) | Out-Null
$rMeta.Collapse(0)
$rMeta.Text = $metaLabel + $metaRest

# Bold just the "Meta description" label portion.
$metaStart = $pMeta.Range.Start
$rBold = $d.Range($metaStart, $metaStart + $metaLabel.Length)
$rBold.Font.Bold = 1

# ---------------------------------------------------------------------------
# 2) Near the end of the document: remove the paragraph that duplicated the
#    title ("Play 1 Million Fortunes Megaways free - Review & Guide", bold),
#    and replace the text of the final (italic) paragraph with the new
#    DALL-E image-generation prompt, keeping its italic formatting.
# ---------------------------------------------------------------------------
$count = $d.Paragraphs.Count
$pDupe = $d.Paragraphs.Item($count - 1)
$pDupe.Range.Delete()

$count2 = $d.Paragraphs.Count
$pLast = $d.Paragraphs.Item($count2)
$rLast = $pLast.Range
$rLast.MoveEnd(1, -1) | Out-Null

$newImagePrompt = 'Create a feature image that fits the game "1 Million Fortunes Megaways" with the following guidelines: - The image should be in cartoon style - The image should feature a happy Maya warrior with glasses DALLE, please create a colorful and engaging feature image for "1 Million Fortunes Megaways" that showcases a happy Maya warrior with glasses. The image should be in a cartoon style that captures the excitement and adventurous nature of the game. Think bold, vibrant colors and eye-catching design elements that will draw in potential players. The Maya warrior should be front and center, surrounded by lotus flowers, gold coins, and other symbols of wealth and prosperity. The image should capture the spirit of the game, conveying the idea of fortune and adventure waiting to be discovered. Be creative and have fun with it!'

$rLast.Text = $newImagePrompt
